$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The latest TestNG run failed for the first three cases: update the
# "Actual Result" and "Status" columns from Pass to Fail for rows 2-4.
$ws.Range("D2:E4").Value = "Fail"
